$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '70.888.08'
$ws.Range("E2").Value = '  +0.09%  '

$ws.Range("D3").Value = '3.848.04'
$ws.Range("E3").Value = '  +1.35%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.00'
$ws.Range("E4").Value = '  -0.04%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '706.85'
$ws.Range("E5").Value = '  +0.79%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '172.53'
$ws.Range("E6").Value = '  -0.15%  '

$ws.Range("D7").Value = '3.845.90'
$ws.Range("E7").Value = '  +1.34%  '

$ws.Range("E9").Value = '  -0.46%  '

$ws.Range("E10").Value = '  -0.30%  '

$ws.Range("E11").Value = '  -1.48%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.0000256'
$ws.Range("E13").Value = '  -0.63%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '36.76'
$ws.Range("E14").Value = '  +0.92%  '

$ws.Range("D15").Value = '4.496.81'
$ws.Range("E15").Value = '  +1.36%  '

$ws.Range("D16").Value = '3.783.23'
$ws.Range("E16").Value = '  -0.26%  '

$ws.Range("D17").Value = '70.943.01'
$ws.Range("E17").Value = '  +0.17%  '

$ws.Range("E18").Value = '  +0.10%  '

$ws.Range("E19").Value = '  +0.95%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '17.33'
$ws.Range("E20").Value = '  -2.93%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '10.66'
$ws.Range("E21").Value = '  -3.86%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '492.95'
$ws.Range("E22").Value = '  +2.21%  '

$ws.Range("E23").Value = '  +0.36%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '85.19'
$ws.Range("E24").Value = '  +0.72%  '

$ws.Range("E25").Value = '  +2.20%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '10.59'
$ws.Range("E26").Value = '  +1.16%  '

$ws.Range("E27").Value = '  -2.07%  '

$ws.Range("E28").Value = '  -3.26%  '

$ws.Range("E29").Value = '  +1.68%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.999'
$ws.Range("E30").Value = '  -0.03%  '

$ws.Range("E31").Value = '  -0.57%  '

$ws.Range("E32").Value = '  -0.63%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '29.45'
$ws.Range("E33").Value = '  -0.21%  '

$ws.Range("E34").Value = '  -1.06%  '

$ws.Range("B35").Value = 'RenzoRestakedETH'
$ws.Range("C35").Value = 'https://coinranking.com/coin/lKlJ_MC5M+renzorestakedeth-ezeth'
$ws.Range("D35").Value = '3.804.13'
$ws.Range("E35").Value = '  +1.54%  '

$ws.Range("B36").Value = 'Aptos'
$ws.Range("C36").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '9.15'
$ws.Range("E36").Value = '  -0.91%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '1.00'
$ws.Range("E37").Value = '  +0.20%  '

$ws.Range("E38").Value = '  +0.26%  '

$ws.Range("E39").Value = '  +6.88%  '

$ws.Range("B40").Value = 'Mantle'
$ws.Range("C40").Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '1.04'
$ws.Range("E40").Value = '  +6.66%  '

$ws.Range("B41").Value = 'Filecoin'
$ws.Range("C41").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '6.05'
$ws.Range("E41").Value = '  +0.19%  '

$ws.Range("E42").Value = '  -3.33%  '

$ws.Range("E44").Value = '  +0.16%  '

$ws.Range("E45").Value = '  -4.17%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '162.52'
$ws.Range("E46").Value = '  +0.18%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '48.71'
$ws.Range("E47").Value = '  -0.64%  '

$ws.Range("E48").Value = '  +1.65%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '415.59'
$ws.Range("E49").Value = '  +1.53%  '

$ws.Range("E50").Value = '  -1.42%  '

$ws.Range("E51").Value = '  +0.74%  '

